$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"15.93194"
$ws.Range("H2").Value = [double]"47.79582"
$ws.Range("I2").Value = [double]"0.9552847657129105"
$ws.Range("J2").Value = [double]"0.9552847657129107"
$ws.Range("M2").Value = [double]"0.9705896666666667"
$ws.Range("N2").Value = [double]"2.911769"
$ws.Range("O2").Value = [double]"0.02073452941466921"
$ws.Range("P2").Value = [double]"0.02073452941466921"
$ws.Range("Q2").Value = [double]"15.46337633395333"
$ws.Range("R2").Value = [double]"139.17038700558"
$ws.Range("S2").Value = [double]"0.01980738007405973"
$ws.Range("T2").Value = [double]"0.01980738007405973"
$ws.Range("G3").Value = [double]"15.93194"
$ws.Range("H3").Value = [double]"47.79582"
$ws.Range("I3").Value = [double]"0.9552847657129105"
$ws.Range("J3").Value = [double]"0.9552847657129107"
$ws.Range("O3").Value = [double]"0.5628689972673966"
$ws.Range("P3").Value = [double]"0.5628689972673966"
$ws.Range("Q3").Value = [double]"419.7758703557999"
$ws.Range("R3").Value = [double]"3777.982833202199"
$ws.Range("S3").Value = [double]"0.5377001781816458"
$ws.Range("T3").Value = [double]"0.537700178181646"
$ws.Range("G4").Value = [double]"15.93194"
$ws.Range("H4").Value = [double]"47.79582"
$ws.Range("I4").Value = [double]"0.9552847657129105"
$ws.Range("J4").Value = [double]"0.9552847657129107"
$ws.Range("M4").Value = [double]"19.49164633333333"
$ws.Range("N4").Value = [double]"58.47493899999999"
$ws.Range("O4").Value = [double]"0.4163964733179342"
$ws.Range("P4").Value = [double]"0.4163964733179341"
$ws.Range("Q4").Value = [double]"310.5397398838866"
$ws.Range("R4").Value = [double]"2794.857658954979"
$ws.Range("S4").Value = [double]"0.397777207457205"
$ws.Range("T4").Value = [double]"0.397777207457205"
$ws.Range("I5").Value = [double]"0.004609931913019111"
$ws.Range("J5").Value = [double]"0.004609931913019112"
$ws.Range("M5").Value = [double]"0.9705896666666667"
$ws.Range("N5").Value = [double]"2.911769"
$ws.Range("O5").Value = [double]"0.02073452941466921"
$ws.Range("P5").Value = [double]"0.02073452941466921"
$ws.Range("Q5").Value = [double]"0.07462184534233333"
$ws.Range("R5").Value = [double]"0.671596608081"
$ws.Range("S5").Value = [double]"9.558476885011707e-05"
$ws.Range("T5").Value = [double]"9.55847688501171e-05"
$ws.Range("I6").Value = [double]"0.004609931913019111"
$ws.Range("J6").Value = [double]"0.004609931913019112"
$ws.Range("O6").Value = [double]"0.5628689972673966"
$ws.Range("P6").Value = [double]"0.5628689972673966"
$ws.Range("S6").Value = [double]"0.002594787753352038"
$ws.Range("T6").Value = [double]"0.002594787753352039"
$ws.Range("I7").Value = [double]"0.004609931913019111"
$ws.Range("J7").Value = [double]"0.004609931913019112"
$ws.Range("M7").Value = [double]"19.49164633333333"
$ws.Range("N7").Value = [double]"58.47493899999999"
$ws.Range("O7").Value = [double]"0.4163964733179342"
$ws.Range("P7").Value = [double]"0.4163964733179341"
$ws.Range("Q7").Value = [double]"1.498576245045666"
$ws.Range("R7").Value = [double]"13.487186205411"
$ws.Range("S7").Value = [double]"0.001919559390816956"
$ws.Range("T7").Value = [double]"0.001919559390816956"
$ws.Range("G8").Value = [double]"0.6688636666666667"
$ws.Range("H8").Value = [double]"2.006591"
$ws.Range("I8").Value = [double]"0.04010530237407027"
$ws.Range("J8").Value = [double]"0.04010530237407027"
$ws.Range("M8").Value = [double]"0.9705896666666667"
$ws.Range("N8").Value = [double]"2.911769"
$ws.Range("O8").Value = [double]"0.02073452941466921"
$ws.Range("P8").Value = [double]"0.02073452941466921"
$ws.Range("Q8").Value = [double]"0.6491921632754445"
$ws.Range("R8").Value = [double]"5.842729469479001"
$ws.Range("S8").Value = [double]"0.0008315645717593631"
$ws.Range("T8").Value = [double]"0.0008315645717593631"
$ws.Range("G9").Value = [double]"0.6688636666666667"
$ws.Range("H9").Value = [double]"2.006591"
$ws.Range("I9").Value = [double]"0.04010530237407027"
$ws.Range("J9").Value = [double]"0.04010530237407027"
$ws.Range("O9").Value = [double]"0.5628689972673966"
$ws.Range("P9").Value = [double]"0.5628689972673966"
$ws.Range("Q9").Value = [double]"17.62326670979"
$ws.Range("R9").Value = [double]"158.60940038811"
$ws.Range("S9").Value = [double]"0.02257403133239867"
$ws.Range("T9").Value = [double]"0.02257403133239867"
$ws.Range("G10").Value = [double]"0.6688636666666667"
$ws.Range("H10").Value = [double]"2.006591"
$ws.Range("I10").Value = [double]"0.04010530237407027"
$ws.Range("J10").Value = [double]"0.04010530237407027"
$ws.Range("M10").Value = [double]"19.49164633333333"
$ws.Range("N10").Value = [double]"58.47493899999999"
$ws.Range("O10").Value = [double]"0.4163964733179342"
$ws.Range("P10").Value = [double]"0.4163964733179341"
$ws.Range("Q10").Value = [double]"13.03725403588322"
$ws.Range("R10").Value = [double]"117.335286322949"
$ws.Range("S10").Value = [double]"0.01669970646991223"
$ws.Range("T10").Value = [double]"0.01669970646991223"
